$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows 14-20 (column B text, column C number).
# Columns D and E repeat the same values used throughout the table
# (adam / krzywo).
$rowsData = @(
    @{ Row = 14; B = "nowa7";  C = 38 },
    @{ Row = 15; B = "nowa8";  C = 39 },
    @{ Row = 16; B = "nowa9";  C = 40 },
    @{ Row = 17; B = "nowa10"; C = 41 },
    @{ Row = 18; B = "nowa11"; C = 42 },
    @{ Row = 19; B = "nowa12"; C = 43 },
    @{ Row = 20; B = "nowa13"; C = 44 }
)

foreach ($item in $rowsData) {
    $r = $item.Row
    $prev = $r - 1

    # Put the value in B so the row gets created, then copy the
    # formatting (only) from the row above (columns B:E) so column A
    # is left untouched (no cell written there), matching the style
    # used by all the other data rows.
    $ws.Range("B$r").Value = $item.B
    $ws.Range("B$prev`:E$prev").Copy()
    $ws.Range("B$r`:E$r").PasteSpecial(-4122)

    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = "adam"
    $ws.Range("E$r").Value = "krzywo"
}

# Row 19 also carries a plain (unstyled) value in column A
$ws.Range("A19").Value = 16

# Update the view: scroll down a bit and move the active selection to A19
$excel.Goto($ws.Range("A7"), $true)
$ws.Range("A19").Select()
